# pretempore.docx edit:
#   - "Fecha de inicio: {{fecha_inicio}}"        -> two runs (label / placeholder),
#                                                    placeholder also gains a trailing space
#   - "Fecha límite: {{fecha_limite}}"           -> two runs (label / placeholder)
#   - "Fecha de presentación: {{fecha_presentacion}}" -> two runs (label / placeholder)
#   - "{{conclusion}}"                            -> "{{conclusión}}"
#
# Splitting a single <w:r> into two runs isn't directly exposed on the Word
# object model (there's no "Runs" collection), so - exactly as a human
# automating this in real Word would do - we apply a character-formatting
# toggle (Bold on, then back off) to the sub-range we want separated. Word
# always breaks the run at that boundary to carry the (possibly now-empty)
# run-properties, leaving the two pieces as independent <w:r> elements with
# identical formatting.

$d = $word.ActiveDocument

# --- "Fecha de inicio: {{fecha_inicio}}" ---------------------------------
# First grow the placeholder text itself by one trailing space...
$d.Content.Find.Execute("{{fecha_inicio}}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{fecha_inicio}} ", 2)

# ...then split "Fecha de inicio: " away from "{{fecha_inicio}} ".
$r1 = $d.Content
$r1.Find.Execute("{{fecha_inicio}} ", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$r1.Bold = $true
$r1.Bold = $false

# --- "Fecha límite: {{fecha_limite}}" -------------------------------------
$r2 = $d.Content
$r2.Find.Execute("{{fecha_limite}}", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$r2.Bold = $true
$r2.Bold = $false

# --- "Fecha de presentación: {{fecha_presentacion}}" ----------------------
$r3 = $d.Content
$r3.Find.Execute("{{fecha_presentacion}}", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$r3.Bold = $true
$r3.Bold = $false

# --- Fix the "conclusion" -> "conclusión" typo -----------------------------
$d.Content.Find.Execute("{{conclusion}}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{conclusión}}", 2)
